$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '27.783.64'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +3.05%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.866.41'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +2.86%  '

$ws.Range("E4").Value = '  +3.40%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '324.57'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +4.27%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '1.037'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +3.15%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.4428'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +3.13%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.07476'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +3.44%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.8860'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +2.64%  '

$ws.Range("E11").Value = '  +2.18%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.890.81'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -12.76%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '5.567'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +3.00%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '6.767'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +2.66%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.07246'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +4.47%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '83.91'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +3.49%  '

$ws.Range("E17").Value = '  +3.00%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '0.000009164'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +3.15%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '1.037'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +3.16%  '

$ws.Range("E20").Value = '  +2.71%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '27.797.63'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +2.95%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '5.323'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +3.02%  '

$ws.Range("E23").Value = '  +3.49%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '1.990'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +5.83%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '158.93'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +3.33%  '

$ws.Range("E26").Value = '  +3.31%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '5.338'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +2.34%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '1.988'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +4.87%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '117.95'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +2.97%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '0.09076'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +1.61%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.7789'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +4.54%  '

$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '3.114'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +10.99%  '

$ws.Range("E33").Value = '  +2.28%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '4.575'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +3.76%  '

$ws.Range("E35").Value = '  +3.30%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '1.158'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +2.49%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.01995'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +3.99%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.05354'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +2.89%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '2.870'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +4.62%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.5206'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +2.22%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.1697'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +2.76%  '

$ws.Range("E42").Value = '  +6.75%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '8.687'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +4.63%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '109.95'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +3.29%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '10.69'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +2.56%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '1.724'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +5.02%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.4718'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +3.46%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '0.06481'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +4.44%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '1.920'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +4.43%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '39.95'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +3.52%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '64.67'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +2.78%  '
